$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Area column (G) - new shared formula group G4:G15, plus standalone G2/G3
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Atotal (H2) sums the Area column
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Summary cells pulling totals together
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Re-enter the existing Q formula range so it becomes a shared formula group (E3:E10)
$ws.Range("E3:E10").Formula = "=(D3-D2)*(B2/100)*C3"

# Update selection to match the authored state
[void]$ws.Range("J2:K2").Select()
